$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: I have recorded that no decision about the movie to be shown on Friday has been made.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies, `"Barbie`" and `"Oppenheimer,`" has been processed successfully.`n"
$ws.Range("D3").Value = "both_movies, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been selected for acquisition.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("D5").Value = "Barbie_was_selected, "
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was reached regarding the movie for Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both movies.`n"
$ws.Range("D7").Value = "both_movies, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be made, as there was no consensus in the discussion.`n"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no agreement. Therefore, no movie will be acquired at this time.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday did not reach a conclusion, so I have recorded the situation as no decision being made.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached, so no further action is taken regarding acquiring movie rights.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision-making process concluded without determining a movie to show.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired for the show on Friday.`n"
$ws.Range("D16").Value = "both_movies, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The rights to the movie `"Barbie`" have been successfully acquired.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded: no movie was selected for Friday.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been confirmed.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has resulted in no agreement.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been confirmed.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision-making committee did not arrive at a definitive choice for Friday's movie, and thus a decision has been logged.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded as a no decision regarding the selection of a movie for Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully made.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been successfully recorded.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision-making process for the movie on Friday resulted in no agreement.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been selected to acquire rights for Friday’s show.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected in the meeting.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The rights for both movies have been acquired.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for showing.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The rights for both movies have been acquired successfully.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision process ended without a conclusive agreement on what movie to show on Friday, so the result is a `"no decision.`"`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The committee did not come to a decision regarding which movie to show on Friday, so there will be no acquisition of movie rights at this time.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been successfully recorded.`n"
$ws.Range("D39").Value = "both_movies, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision for the movie to be shown on Friday has not been made, resulting in no agreement.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" have been acquired for the screening on Friday.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: I have successfully acquired the rights for both movies, `"Oppenheimer`" and `"Barbie,`" to be shown on Friday.`n"
$ws.Range("D44").Value = "both_movies, "
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both movies.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: I have recorded the decision as no decision was made regarding which movie to show on Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday ended without a clear selection. Therefore, the appropriate outcome is recorded as no decision.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" has been recorded.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: I have successfully called the function to acquire the rights for both movies.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been registered, indicating that no movie will be acquired at this time.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired for the event.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision,`" indicating that no movie was selected for Friday.`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been successfully selected for acquisition.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" has been recorded.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been made.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" has been recorded.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision was recorded successfully, indicating that no movie was selected for Friday.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: No decision was made regarding the movie selection.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
